# Update "想去人数" (number of people interested) figures that changed
# between the two scrapes of the 南宁-漫展信息 data, per the commit
# "Update gh-pages to output generated at 456a3b4".
#
# Sheet "展览"   (index 1): F2 958->962, F3 1877->1905, F4 421->424
# Sheet "全部类型" (index 4): F4 958->962, F5 1877->1905, F6 421->424

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 962
$wsExhibit.Range("F3").Value = 1905
$wsExhibit.Range("F4").Value = 424

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 962
$wsAll.Range("F5").Value = 1905
$wsAll.Range("F6").Value = 424
